$d = $word.ActiveDocument

$newBodyXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>Modificacion</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve"> por Andrey Palma</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>Modificacion</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve"> por Andrey Palma</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>Modificacion</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve"> por Andrey Palma</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>Modificacion</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve"> por Andrey Palma</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>Modificacion</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve"> por Andrey Palma</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>conservantes. Ahora puedes preparar tus propios batidos nutritivos en casa de una manera muy fácil. Sólo necesita</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>s una batidora de alta potencia</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>, y estás listo para comenzar.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve">con la almacenada en el elemento seguro. El iluminador infrarrojo ayuda a identificar la </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve">cara del usuario en la oscuridad. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve">El proyector de puntos sitúa más de 30 000 puntos invisibles sobre la cara para crear un mapa facial único. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t xml:space="preserve">El sistema no funcionará con los ojos cerrados o no mirando al teléfono, impidiendo el acceso si no está mirándose </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="es-CR"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
          <w:lang w:val="es-CR"/>
        </w:rPr>
        <w:t>activamente la pantalla del dispositivo.</w:t>
      </w:r>
    </w:p>
    
'@

$d.Content.InsertXML($newBodyXml)
